$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F ("想去人数" / "want to go" count) updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 173
$wsExpo.Range("F4").Value  = 176
$wsExpo.Range("F5").Value  = 4942
$wsExpo.Range("F8").Value  = 7
$wsExpo.Range("F10").Value = 502
$wsExpo.Range("F14").Value = 3526
$wsExpo.Range("F16").Value = 130
$wsExpo.Range("F18").Value = 76
$wsExpo.Range("F19").Value = 2591
$wsExpo.Range("F27").Value = 258

# Sheet "全部类型" (All Types) - same underlying records, column F updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 173
$wsAll.Range("F4").Value  = 176
$wsAll.Range("F6").Value  = 4942
$wsAll.Range("F9").Value  = 7
$wsAll.Range("F11").Value = 502
$wsAll.Range("F15").Value = 3526
$wsAll.Range("F17").Value = 130
$wsAll.Range("F19").Value = 76
$wsAll.Range("F20").Value = 2591
$wsAll.Range("F25").Value = 43
$wsAll.Range("F28").Value = 258
